# Generate Report for Handoff
#
# Refresh the "Handback DateTime" timestamp column on each sheet for the
# rows whose handback re-ran: row 7 and rows 10-16 (row 8 "In Translation"
# and row 9 "In Translation" keep their own distinct timestamps and are
# left untouched).
#
# Overview sheet uses column D; zh-cn / de-de sheets use column E.
# zh-cn gets its own distinct refreshed timestamp (08:26:42); Overview and
# de-de share the same refreshed timestamp (08:26:46).

$wb = $excel.ActiveWorkbook

$rows = @(7, 10, 11, 12, 13, 14, 15, 16)

$ws = $wb.Worksheets.Item("Overview")
foreach ($r in $rows) {
    $ws.Range("D$r").Value = "2016-03-21 08:26:46"
}

$ws = $wb.Worksheets.Item("zh-cn")
foreach ($r in $rows) {
    $ws.Range("E$r").Value = "2016-03-21 08:26:42"
}

$ws = $wb.Worksheets.Item("de-de")
foreach ($r in $rows) {
    $ws.Range("E$r").Value = "2016-03-21 08:26:46"
}
